$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")
$ws.Activate()

$data = @(
    ,@(2, "4398", "3", "", "", "", "NO")
    ,@(3, "4399", "4", "4", "0", "15.13%", "NO")
    ,@(4, "4400", "4", "5", "0", "18.88%", "NO")
    ,@(5, "4415", "4", "2", "0", "18.89%", "NO")
    ,@(6, "4419", "", "", "", "", "NO")
    ,@(7, "4421", "", "", "", "", "NO")
    ,@(8, "4423", "5", "2", "0", "21.71%", "NO")
    ,@(9, "4429", "4", "1", "0", "7.14%", "NO")
    ,@(10, "4430", "4", "3", "0", "23.19%", "NO")
    ,@(11, "4431", "4", "1", "0", "6.56%", "NO")
    ,@(12, "4435", "", "", "", "", "NO")
    ,@(13, "4436", "", "", "", "", "NO")
    ,@(14, "4437", "1", "1", "0", "2.42%", "NO")
    ,@(15, "4564", "4", "1", "0", "7.99%", "NO")
    ,@(16, "4565", "", "", "", "", "NO")
    ,@(17, "4567", "4", "1", "0", "1.90%", "NO")
    ,@(18, "4594", "", "", "", "", "NO")
    ,@(19, "4597", "5", "0", "0", "9.52%", "NO")
    ,@(20, "4600", "", "", "", "", "NO")
    ,@(21, "4601", "4", "0", "0", "5.51%", "NO")
    ,@(22, "4603", "5", "2", "0", "18.90%", "NO")
    ,@(23, "4647", "", "", "", "", "NO")
    ,@(24, "4648", "4", "0", "0", "2.56%", "NO")
    ,@(25, "4649", "4", "2", "0", "19.48%", "NO")
    ,@(26, "4660", "", "", "", "", "NO")
    ,@(27, "4663", "", "", "", "", "NO")
    ,@(28, "4666", "", "", "", "", "NO")
    ,@(29, "4725", "4", "1", "0", "7.98%", "NO")
    ,@(30, "4728", "4", "", "", "", "NO")
    ,@(31, "4732", "5", "1", "1", "10.41%", "NO")
)

foreach ($row in $data) {
    $r = $row[0]
    $matchCode = $row[1]
    $battingPos = $row[2]
    $num4 = $row[3]
    $num6 = $row[4]
    $pct = $row[5]
    $mom = $row[6]

    # Column A - MATCH_CODE, always text even though numeric-looking
    $ws.Cells.Item($r, 1).Value = "'" + $matchCode

    # Column B - BATTING_POSITION, numeric when present, text-empty otherwise
    if ($battingPos -eq "") {
        $ws.Cells.Item($r, 2).Value = "'"
    } else {
        $ws.Cells.Item($r, 2).Value = [int]$battingPos
    }

    # Column C - NUM_4, text
    if ($num4 -eq "") {
        $ws.Cells.Item($r, 3).Value = "'"
    } else {
        $ws.Cells.Item($r, 3).Value = "'" + $num4
    }

    # Column D - NUM_6, text
    if ($num6 -eq "") {
        $ws.Cells.Item($r, 4).Value = "'"
    } else {
        $ws.Cells.Item($r, 4).Value = "'" + $num6
    }

    # Column E - PERCENT_RUNS_OF_TOTAL, text
    if ($pct -eq "") {
        $ws.Cells.Item($r, 5).Value = "'"
    } else {
        $ws.Cells.Item($r, 5).Value = "'" + $pct
    }

    # Column F - MAN_OF_MATCH, plain text
    $ws.Cells.Item($r, 6).Value = $mom
}

Write-Host "Edit complete"
